$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set transponder_id (column A) and name (column B) values for rows 2-45
$ws.Cells.Item(2, 1).Value = "FF30660"
$ws.Cells.Item(2, 2).Value = "Jesus"
$ws.Cells.Item(3, 1).Value = "GN87210"
$ws.Cells.Item(3, 2).Value = "Brian"
$ws.Cells.Item(4, 1).Value = "GR11874"
$ws.Cells.Item(4, 2).Value = "Mark"
$ws.Cells.Item(5, 1).Value = "GT49254"
$ws.Cells.Item(5, 2).Value = "Emily"
$ws.Cells.Item(6, 1).Value = "PX83485"
$ws.Cells.Item(6, 2).Value = "David"
$ws.Cells.Item(7, 1).Value = "RR01344"
$ws.Cells.Item(7, 2).Value = "Charles"
$ws.Cells.Item(8, 1).Value = "RZ04998"
$ws.Cells.Item(8, 2).Value = "Megan"
$ws.Cells.Item(9, 1).Value = "RZ13738"
$ws.Cells.Item(9, 2).Value = "Maria"
$ws.Cells.Item(10, 1).Value = "SC98900"
$ws.Cells.Item(10, 2).Value = "Charles"
$ws.Cells.Item(11, 1).Value = "SF16020"
$ws.Cells.Item(11, 2).Value = "Jessica"
$ws.Cells.Item(12, 1).Value = "CS34132"
$ws.Cells.Item(12, 2).Value = "Tracey"
$ws.Cells.Item(13, 1).Value = "CR48188"
$ws.Cells.Item(13, 2).Value = "Andrea"
$ws.Cells.Item(14, 1).Value = "RZ10761"
$ws.Cells.Item(14, 2).Value = "Heather"
$ws.Cells.Item(15, 1).Value = "KT68858"
$ws.Cells.Item(15, 2).Value = "Jason"
$ws.Cells.Item(16, 1).Value = "HC20410"
$ws.Cells.Item(16, 2).Value = "Michael"
$ws.Cells.Item(17, 1).Value = "NS43328"
$ws.Cells.Item(17, 2).Value = "Tyler"
$ws.Cells.Item(18, 1).Value = "RR15537"
$ws.Cells.Item(18, 2).Value = "Amanda"
$ws.Cells.Item(19, 1).Value = "RZ43503"
$ws.Cells.Item(19, 2).Value = "James"
$ws.Cells.Item(20, 1).Value = "LZ73964"
$ws.Cells.Item(20, 2).Value = "Kristin"
$ws.Cells.Item(21, 1).Value = "KP55121"
$ws.Cells.Item(21, 2).Value = "Lisa"
$ws.Cells.Item(22, 1).Value = "CW09719"
$ws.Cells.Item(22, 2).Value = "Michelle"
$ws.Cells.Item(23, 1).Value = "SF04172"
$ws.Cells.Item(23, 2).Value = "Joshua"
$ws.Cells.Item(24, 1).Value = "RR10973"
$ws.Cells.Item(24, 2).Value = "James"
$ws.Cells.Item(25, 1).Value = "HH54287"
$ws.Cells.Item(25, 2).Value = "Patricia"
$ws.Cells.Item(26, 1).Value = "RR18114"
$ws.Cells.Item(26, 2).Value = "Debbie"
$ws.Cells.Item(27, 1).Value = "KP31430"
$ws.Cells.Item(27, 2).Value = "Kathleen"
$ws.Cells.Item(28, 1).Value = "PF79070"
$ws.Cells.Item(28, 2).Value = "Brandy"
$ws.Cells.Item(29, 1).Value = "RF48933"
$ws.Cells.Item(29, 2).Value = "Debbie"
$ws.Cells.Item(30, 1).Value = "LR40395"
$ws.Cells.Item(30, 2).Value = "Shane"
$ws.Cells.Item(31, 1).Value = "RZ02081"
$ws.Cells.Item(31, 2).Value = "Robert"
$ws.Cells.Item(32, 1).Value = "GK67167"
$ws.Cells.Item(32, 2).Value = "George"
$ws.Cells.Item(33, 1).Value = "CX19285"
$ws.Cells.Item(33, 2).Value = "Tracie"
$ws.Cells.Item(34, 1).Value = "SC75193"
$ws.Cells.Item(34, 2).Value = "Christopher"
$ws.Cells.Item(35, 1).Value = "RZ39892"
$ws.Cells.Item(35, 2).Value = "Richard"
$ws.Cells.Item(36, 1).Value = "RX72266"
$ws.Cells.Item(36, 2).Value = "Victor"
$ws.Cells.Item(37, 1).Value = "FW08700"
$ws.Cells.Item(37, 2).Value = "Kevin"
$ws.Cells.Item(38, 1).Value = "VG76471"
$ws.Cells.Item(38, 2).Value = "Cody"
$ws.Cells.Item(39, 1).Value = "KR51169"
$ws.Cells.Item(39, 2).Value = "William"
$ws.Cells.Item(40, 1).Value = "NK39779"
$ws.Cells.Item(40, 2).Value = "Douglas"
$ws.Cells.Item(41, 1).Value = "RZ41573"
$ws.Cells.Item(41, 2).Value = "Jennifer"
$ws.Cells.Item(42, 1).Value = "RX45371"
$ws.Cells.Item(42, 2).Value = "Lisa"
$ws.Cells.Item(43, 1).Value = "SC90979"
$ws.Cells.Item(43, 2).Value = "Timothy"
$ws.Cells.Item(44, 1).Value = "GZ52385"
$ws.Cells.Item(44, 2).Value = "Sean"
$ws.Cells.Item(45, 1).Value = "SF17035"
$ws.Cells.Item(45, 2).Value = "Catherine"

# Remove now-unused rows 46-48 (sheet shrank from 48 to 45 rows)
$ws.Range("A46:B48").ClearContents()
